$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2
Set-TextValue 'D2' '44.171.03'
Set-TextValue 'E2' '  +1.00%  '

# Row 3
Set-TextValue 'D3' '2.368.21'
Set-TextValue 'E3' '  +0.76%  '

# Row 4
Set-TextValue 'E4' '  +0.15%  '

# Row 5
Set-TextValue 'D5' '0.700'
Set-TextValue 'E5' '  +6.54%  '

# Row 6
Set-TextValue 'D6' '241.90'
Set-TextValue 'E6' '  +3.08%  '

# Row 7
Set-TextValue 'D7' '76.90'
Set-TextValue 'E7' '  +4.35%  '

# Row 8
Set-TextValue 'E8' '  +0.00%  '

# Row 9
Set-TextValue 'E9' '  +16.67%  '

# Row 10
Set-TextValue 'E10' '  +4.15%  '

# Row 11
Set-TextValue 'D11' '57.39'

# Row 12
Set-TextValue 'D12' '33.61'
Set-TextValue 'E12' '  +21.73%  '

# Row 13
Set-TextValue 'D13' '7.50'
Set-TextValue 'E13' '  +12.68%  '

# Row 14
Set-TextValue 'D14' '0.109'
Set-TextValue 'E14' '  +2.25%  '

# Row 15
Set-TextValue 'D15' '2.720.54'
Set-TextValue 'E15' '  +0.95%  '

# Row 16
Set-TextValue 'D16' '16.78'
Set-TextValue 'E16' '  -1.06%  '

# Row 17
Set-TextValue 'D17' '0.933'
Set-TextValue 'E17' '  +5.20%  '

# Row 18
Set-TextValue 'D18' '2.358.64'
Set-TextValue 'E18' '  +0.50%  '

# Row 19
Set-TextValue 'D19' '44.055.43'
Set-TextValue 'E19' '  +1.09%  '

# Row 21
Set-TextValue 'D21' '6.73'
Set-TextValue 'E21' '  +6.58%  '

# Row 22
Set-TextValue 'D22' '77.99'
Set-TextValue 'E22' '  +2.49%  '

# Row 23
Set-TextValue 'D23' '260.61'
Set-TextValue 'E23' '  +4.06%  '

# Row 24
Set-TextValue 'E24' '  +0.17%  '

# Row 25
Set-TextValue 'B25' 'WEMIXToken'
Set-TextValue 'C25' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D25' '3.74'
Set-TextValue 'E25' '  -2.13%  '

# Row 26
Set-TextValue 'B26' 'PancakeSwap'
Set-TextValue 'C26' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D26' '2.53'
Set-TextValue 'E26' '  +2.55%  '

# Row 27
Set-TextValue 'D27' '1.80'
Set-TextValue 'E27' '  +16.99%  '

# Row 28
Set-TextValue 'D28' '10.96'
Set-TextValue 'E28' '  +6.94%  '

# Row 29
Set-TextValue 'B29' 'EthereumClassic'
Set-TextValue 'C29' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D29' '23.27'
Set-TextValue 'E29' '  +3.85%  '

# Row 30
Set-TextValue 'B30' 'Toncoin'
Set-TextValue 'C30' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D30' '2.23'
Set-TextValue 'E30' '  -0.53%  '

# Row 31
Set-TextValue 'D31' '175.39'
Set-TextValue 'E31' '  +1.74%  '

# Row 32
Set-TextValue 'E32' '  -3.36%  '

# Row 33
Set-TextValue 'E33' '  +5.41%  '

# Row 34
Set-TextValue 'E34' '  +6.02%  '

# Row 35
Set-TextValue 'D35' '0.0765'
Set-TextValue 'E35' '  +9.52%  '

# Row 36
Set-TextValue 'E36' '  +7.36%  '

# Row 37
Set-TextValue 'D37' '3.81'
Set-TextValue 'E37' '  +1.39%  '

# Row 38
Set-TextValue 'E38' '  +0.74%  '

# Row 39
Set-TextValue 'D39' '6.45'
Set-TextValue 'E39' '  +1.29%  '

# Row 40
Set-TextValue 'E40' '  +7.39%  '

# Row 41
Set-TextValue 'D41' '0.223'
Set-TextValue 'E41' '  +23.77%  '

# Row 42
Set-TextValue 'B42' 'InjectiveProtocol'
Set-TextValue 'C42' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D42' '19.45'
Set-TextValue 'E42' '  -0.19%  '

# Row 43
Set-TextValue 'B43' 'FraxShare'
Set-TextValue 'C43' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D43' '9.19'
Set-TextValue 'E43' '  +3.43%  '

# Row 44
Set-TextValue 'B44' 'Cronos'
Set-TextValue 'C44' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D44' '0.108'
Set-TextValue 'E44' '  +12.03%  '

# Row 45
Set-TextValue 'B45' 'FTXToken'
Set-TextValue 'C45' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue 'D45' '4.91'
Set-TextValue 'E45' '  +10.67%  '

# Row 46
Set-TextValue 'B46' 'BinanceUSD'
Set-TextValue 'C46' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D46' '1.00'
Set-TextValue 'E46' '  -0.04%  '

# Row 47
Set-TextValue 'D47' '2.54'
Set-TextValue 'E47' '  +11.10%  '

# Row 48
Set-TextValue 'D48' '1.27'
Set-TextValue 'E48' '  +4.45%  '

# Row 49
Set-TextValue 'B49' 'ARBITRUM'
Set-TextValue 'C49' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D49' '1.20'
Set-TextValue 'E49' '  +3.19%  '

# Row 50
Set-TextValue 'B50' 'Aave'
Set-TextValue 'C50' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D50' '102.51'
Set-TextValue 'E50' '  +3.90%  '

# Row 51
Set-TextValue 'D51' '56.27'
Set-TextValue 'E51' '  +9.57%  '
